$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect("D382")

# Update the confidential disclaimer date (2021-04-29 -> 2021-04-30)
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-04-30 for illustrative purposes only and are subject to change."

# Update Weight (D) / Percent Change (E) values for rows 2-56
$ws.Range("D2").Value = 0.01542824277730473
$ws.Range("E2").Value = 0.003723088544757891
$ws.Range("D3").Value = 0.05351436520169796
$ws.Range("E3").Value = -0.0011206144078173
$ws.Range("D4").Value = 0.01431824645797495
$ws.Range("E4").Value = 0.01257563685695229
$ws.Range("D5").Value = 0.009460605833814324
$ws.Range("E5").Value = 0.009125277017338007
$ws.Range("D6").Value = 0.01508924081816385
$ws.Range("E6").Value = 0.001528414617201435
$ws.Range("D7").Value = 0.01920231042720004
$ws.Range("E7").Value = 0.006867406233491913
$ws.Range("D8").Value = 0.004094446854755055
$ws.Range("E8").Value = -0.037320120183438
$ws.Range("D9").Value = 0.006542136580113606
$ws.Range("E9").Value = -0.02884289107567028
$ws.Range("D10").Value = 0.01424918193878663
$ws.Range("E10").Value = -0.006075949367088551
$ws.Range("D11").Value = 0.008178595696381135
$ws.Range("E11").Value = 0.005202442886224734
$ws.Range("D12").Value = 0.01506223165798127
$ws.Range("E12").Value = -0.004175878162613555
$ws.Range("D13").Value = 0.002936321198615915
$ws.Range("E13").Value = 0.02756339581036382
$ws.Range("D14").Value = 0.00625064731028933
$ws.Range("E14").Value = -0.01465002712967989
$ws.Range("D15").Value = 0.01387614103804579
$ws.Range("E15").Value = -0.008892325536439261
$ws.Range("D16").Value = 0.01024498144459596
$ws.Range("E16").Value = -0.03521126760563398
$ws.Range("D17").Value = 0.02149309419802763
$ws.Range("E17").Value = -0.001742947517913707
$ws.Range("D18").Value = 0.008144001772037701
$ws.Range("E18").Value = -0.007473366194943654
$ws.Range("D19").Value = 0.01620056292960305
$ws.Range("E19").Value = -0.008952496954933054
$ws.Range("D20").Value = 0.01123531731795705
$ws.Range("E20").Value = 0.00312843029637766
$ws.Range("D21").Value = 0.007150520996463574
$ws.Range("E21").Value = 0.002233567326100916
$ws.Range("D22").Value = 0.01318614332602998
$ws.Range("E22").Value = 0.00039282440749
$ws.Range("D23").Value = 0.01921100515684785
$ws.Range("E23").Value = -0.001343328855570425
$ws.Range("D24").Value = 0.00961109865442066
$ws.Range("E24").Value = -0.03488696622941667
$ws.Range("D25").Value = 0.02066151422167576
$ws.Range("E25").Value = -0.006130208737487264
$ws.Range("D26").Value = 0.01319045985847925
$ws.Range("E26").Value = 0.0006872177498529197
$ws.Range("D27").Value = 0.02173888988863892
$ws.Range("E27").Value = -0.03201397881336088
$ws.Range("D28").Value = 0.0576170751328559
$ws.Range("E28").Value = -0.01513335331135723
$ws.Range("D29").Value = 0.02102043389483256
$ws.Range("E29").Value = -0.02081949058693244
$ws.Range("D30").Value = 0.03086197371729497
$ws.Range("E30").Value = -0.02893222506393867
$ws.Range("D31").Value = 0.01627690388920585
$ws.Range("E31").Value = -0.0271937202130641
$ws.Range("D32").Value = 0.01343039739790939
$ws.Range("E32").Value = -0.001221320771176759
$ws.Range("D33").Value = 0.01954871215625396
$ws.Range("E33").Value = -0.02259822848058157
$ws.Range("D34").Value = 0.04426468364277975
$ws.Range("E34").Value = -0.01640783028803572
$ws.Range("D35").Value = 0.01082956326772567
$ws.Range("E35").Value = -0.004783054321831237
$ws.Range("D36").Value = 0.009573205665990996
$ws.Range("E36").Value = 0.005500946884299696
$ws.Range("D37").Value = 0.01072504151770406
$ws.Range("E37").Value = 0.01164294954721856
$ws.Range("D38").Value = 0.007150027678469372
$ws.Range("E38").Value = 0.001940491591203131
$ws.Range("D39").Value = 0.01178792596858752
$ws.Range("E39").Value = -0.01410848941863296
$ws.Range("D40").Value = 0.01681832038784286
$ws.Range("E40").Value = -0.01068791294209093
$ws.Range("D41").Value = 0.01699098168581366
$ws.Range("E41").Value = -0.01441180526823893
$ws.Range("D42").Value = 0.03373647600476241
$ws.Range("E42").Value = -0.008609108107367369
$ws.Range("D43").Value = 0.01122042158112211
$ws.Range("E43").Value = -0.007910817437212558
$ws.Range("D44").Value = 0.02263916439568345
$ws.Range("E44").Value = -0.01393228067212715
$ws.Range("D45").Value = 0.01293914517280783
$ws.Range("E45").Value = -0.02163889425035936
$ws.Range("D46").Value = 0.008215409551698482
$ws.Range("E46").Value = -0.01942173883672227
$ws.Range("D47").Value = 0.01262419246588431
$ws.Range("E47").Value = 0.002813557831802038
$ws.Range("D48").Value = 0.009706648183422715
$ws.Range("E48").Value = -0.01037735849056587
$ws.Range("D49").Value = 0.0153049441111288
$ws.Range("E49").Value = -0.00969395155441688
$ws.Range("D50").Value = 0.008521112546230711
$ws.Range("E50").Value = -0.01070308174940027
$ws.Range("D51").Value = 0.01188853200703015
$ws.Range("E51").Value = -0.03133931211195418
$ws.Range("D52").Value = 0.008527063194535776
$ws.Range("E52").Value = -0.01038103571712679
$ws.Range("D53").Value = 0.009661262927956105
$ws.Range("E53").Value = -0.005361451166115527
$ws.Range("D54").Value = 0.133942001900848
$ws.Range("E54").Value = 0.00009852216748762466
$ws.Range("D55").Value = 0.04390807639772077
$ws.Range("E55").Value = -0.007043094186462873
$ws.Range("E56").Value = -0.008414891315719264

$ws.Protect("D382")
